$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a 20x20 distance matrix in B2:U21 (row/col headers are city
# names). Truncate every numeric distance value to an integer (drop the
# decimal portion) - e.g. 1045.501704723878 -> 1045.
for ($r = 2; $r -le 21; $r++) {
    for ($c = 2; $c -le 21; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        $cell.Value2 = [Math]::Floor($val)
    }
}
